$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 77, shifting existing rows 77:175 down to 78:176.
$ws.Rows("77:77").Insert()

# Populate the newly inserted row 77 with the new data record.
$ws.Range("A77").Value = 4
$ws.Range("B77").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C77").Value = "Los Lagos"
$ws.Range("D77").Value = 45174
$ws.Range("E77").Value = 10
$ws.Range("F77").Value = 100112022
$ws.Range("G77").Value = "Arveja Verde"
$ws.Range("H77").Value = "Perfection"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 35
$ws.Range("K77").Value = 37000
$ws.Range("L77").Value = 37000
$ws.Range("M77").Value = 37000
$ws.Range("N77").Value = "`$/malla 25 kilos"
$ws.Range("O77").Value = "Provincia de Limarí"
$ws.Range("P77").Value = 1480
$ws.Range("Q77").Value = 25
$ws.Range("R77").Value = "Hortaliza"
